$wb = $excel.ActiveWorkbook

# --- Rename sheets (per commit diff) ---
$wb.Worksheets.Item("Include ValueSets").Name = "Include ValueSet #0"
$wb.Worksheets.Item("Exclude from Consent Category").Name = "Exclude #1"
$wb.Worksheets.Item("Exclude from LOINC").Name = "Exclude #2"

# --- Metadata sheet updates ---
$ws = $wb.Worksheets.Item("Metadata")

# Simple value updates (these cells already exist, so their existing
# style ("s" attribute) is preserved automatically).
$ws.Range("B3").Value = "0.2.2"
$ws.Range("B8").Value = "2024-09-11T16:17:59-05:00"
$ws.Range("B10").Value = "MITRE (https://www.mitre.org)"

# Insert a new row at position 11 for "Jurisdiction" - this shifts the
# old rows 11-14 (Description/Purpose/Copyright/Immutable) down to 12-15,
# keeping their original formatting intact since they are simply moved.
$ws.Rows.Item(11).Insert()

# The freshly inserted row 11 has no formatting yet; copy the format
# (only) from the row right below it (now row 12, which kept its
# original style) so the new row matches the rest of the table instead
# of minting a brand-new style entry.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
